$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# BDbDT: extend the model from 2050 through 2070 (columns AG:AZ, 20 more
# years) for rows 1-9, using FORECAST.ETS for the two computed population
# rows (2 & 3) and flat zeros for the remaining (no-data) demographic rows.
# ---------------------------------------------------------------------------
$bdbdt = $wb.Worksheets.Item("BDbDT")

# Row 1: year headers 2051..2070
$year = 2051
for ($col = 33; $col -le 52; $col++) {
    $bdbdt.Cells.Item(1, $col).Value = $year
    $year = $year + 1
}

# Rows 2 & 3: FORECAST.ETS trend continuation, anchored on the fixed
# 2020-2050 history ($B$<row>:$AF$<row>) and the fixed year axis ($B$1:$AF$1)
for ($col = 33; $col -le 52; $col++) {
    $colLetter = $bdbdt.Cells.Item(1, $col).Address($false, $false, 1, $false)
    $bdbdt.Cells.Item(2, $col).Formula = "=FORECAST.ETS(" + $colLetter + ",`$B`$2:`$AF`$2,`$B`$1:`$AF`$1)"
    $bdbdt.Cells.Item(3, $col).Formula = "=FORECAST.ETS(" + $colLetter + ",`$B`$3:`$AF`$3,`$B`$1:`$AF`$1)"
}

# Rows 4-9: no statistics by race/ethnicity available -> continue the flat
# zero series, carrying forward the same number format as column AF.
for ($row = 4; $row -le 9; $row++) {
    $bdbdt.Range($bdbdt.Cells.Item($row, 32), $bdbdt.Cells.Item($row, 32)).Copy()
    $bdbdt.Range($bdbdt.Cells.Item($row, 33), $bdbdt.Cells.Item($row, 52)).PasteSpecial(-4122)
    for ($col = 33; $col -le 52; $col++) {
        $bdbdt.Cells.Item($row, $col).Value = 0
    }
}
$bdbdt.Application.CutCopyMode = 0

# ---------------------------------------------------------------------------
# View bookkeeping: restore "About"/"IBGE" scroll position and make sure
# BDbDT stays the active sheet/tab with the new selection in view, matching
# the saved state of the edited workbook.
# ---------------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")
$about.Range("A16:C16").Select()

$ibge = $wb.Worksheets.Item("IBGE")
$ibge.Range("G8").Select()

$bdbdt.Activate()
$bdbdt.Range("AP2").Select()
